$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.417.29'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.83%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.669.38'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.14%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.32'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.09%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3965'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.92%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3929'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.45'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +6.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.396'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.001'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08569'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.79%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.55'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.311'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.943'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +7.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001336'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +5.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.665.19'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.04'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07031'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.64'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.997'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.0000'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.429.82'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.475'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +6.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.070'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +15.70%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.56'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.46'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '142.87'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.444'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.024'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -6.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.543'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.844.31'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.068'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +14.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.03098'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +8.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08302'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +5.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.924'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '11.17'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +13.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2769'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +4.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09269'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7722'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.78'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +7.27%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.70'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7130'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +5.06%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.125'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.22%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08442'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '137.07'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.44%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.97%  '
